# Replace Greek-letter abbreviations in biomarker labels with plain-English
# (ASCII) equivalents, per commit "greek letters in english".

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "Vit D (μmol/L)";  New = "Vit D (umol/L)" },
    @{ Old = "Ferritin (μg/L)"; New = "Ferritin (ug/L)" },
    @{ Old = "IL-1β (pg/ml)";   New = "IL-1B (pg/ml)" },
    @{ Old = "TNF-α (pg/ml)";   New = "TNF-a (pg/ml)" },
    @{ Old = "IFN-γ (pg/ml)";   New = "IFN-y (pg/ml)" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.New, 2)
}
